$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TEST_CASES")
$ws.Columns("X:Z").Insert()
$ws.Range("X1").Value = "TC_KIND"
$ws.Range("Y1").Value = "TC_SCRIPTING_LANGUAGE"
$ws.Range("Z1").Value = "TC_SCRIPT"
